$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 1141 (shifts old 1141-1199 down to 1145-1203)
$ws.Range("A1141:A1144").EntireRow.Insert()

# Row 1141
$ws.Cells.Item(1141, 1).Value = 3
$ws.Cells.Item(1141, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1141, 3).Value = "Coquimbo"
$ws.Cells.Item(1141, 4).Value = 44753
$ws.Cells.Item(1141, 5).Value = 5
$ws.Cells.Item(1141, 6).Value = 100112020
$ws.Cells.Item(1141, 7).Value = "Tomate"
$ws.Cells.Item(1141, 8).Value = "Larga vida"
$ws.Cells.Item(1141, 9).Value = "Primera"
$ws.Cells.Item(1141, 10).Value = 510
$ws.Cells.Item(1141, 11).Value = 10000
$ws.Cells.Item(1141, 12).Value = 11000
$ws.Cells.Item(1141, 13).Value = 10510
$ws.Cells.Item(1141, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(1141, 15).Value = "Limache"
$ws.Cells.Item(1141, 16).Value = 584
$ws.Cells.Item(1141, 17).Value = 18
$ws.Cells.Item(1141, 18).Value = "Hortaliza"

# Row 1142
$ws.Cells.Item(1142, 1).Value = 3
$ws.Cells.Item(1142, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1142, 3).Value = "Coquimbo"
$ws.Cells.Item(1142, 4).Value = 44753
$ws.Cells.Item(1142, 5).Value = 5
$ws.Cells.Item(1142, 6).Value = 100112020
$ws.Cells.Item(1142, 7).Value = "Tomate"
$ws.Cells.Item(1142, 8).Value = "Larga vida"
$ws.Cells.Item(1142, 9).Value = "Primera"
$ws.Cells.Item(1142, 10).Value = 380
$ws.Cells.Item(1142, 11).Value = 5000
$ws.Cells.Item(1142, 12).Value = 5000
$ws.Cells.Item(1142, 13).Value = 5000
$ws.Cells.Item(1142, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(1142, 15).Value = "Limache"
$ws.Cells.Item(1142, 16).Value = 417
$ws.Cells.Item(1142, 17).Value = 12
$ws.Cells.Item(1142, 18).Value = "Hortaliza"

# Row 1143
$ws.Cells.Item(1143, 1).Value = 3
$ws.Cells.Item(1143, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1143, 3).Value = "Coquimbo"
$ws.Cells.Item(1143, 4).Value = 44753
$ws.Cells.Item(1143, 5).Value = 5
$ws.Cells.Item(1143, 6).Value = 100112020
$ws.Cells.Item(1143, 7).Value = "Tomate"
$ws.Cells.Item(1143, 8).Value = "Larga vida"
$ws.Cells.Item(1143, 9).Value = "Segunda"
$ws.Cells.Item(1143, 10).Value = 468
$ws.Cells.Item(1143, 11).Value = 7500
$ws.Cells.Item(1143, 12).Value = 8000
$ws.Cells.Item(1143, 13).Value = 7746
$ws.Cells.Item(1143, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(1143, 15).Value = "Limache"
$ws.Cells.Item(1143, 16).Value = 430
$ws.Cells.Item(1143, 17).Value = 18
$ws.Cells.Item(1143, 18).Value = "Hortaliza"

# Row 1144
$ws.Cells.Item(1144, 1).Value = 3
$ws.Cells.Item(1144, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1144, 3).Value = "Coquimbo"
$ws.Cells.Item(1144, 4).Value = 44753
$ws.Cells.Item(1144, 5).Value = 5
$ws.Cells.Item(1144, 6).Value = 100112020
$ws.Cells.Item(1144, 7).Value = "Tomate"
$ws.Cells.Item(1144, 8).Value = "Larga vida"
$ws.Cells.Item(1144, 9).Value = "Tercera"
$ws.Cells.Item(1144, 10).Value = 140
$ws.Cells.Item(1144, 11).Value = 5000
$ws.Cells.Item(1144, 12).Value = 5000
$ws.Cells.Item(1144, 13).Value = 5000
$ws.Cells.Item(1144, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(1144, 15).Value = "Limache"
$ws.Cells.Item(1144, 16).Value = 278
$ws.Cells.Item(1144, 17).Value = 18
$ws.Cells.Item(1144, 18).Value = "Hortaliza"
